$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B-E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values for columns B-E
$ws.Range("B2").Value = 83.685038754406875
$ws.Range("C2").Value = 55.699168891699358
$ws.Range("D2").Value = 52.557941815545774
$ws.Range("E2").Value = 51.337072624918413

# Update row 3 values for columns B-E
$ws.Range("B3").Value = 72.869242233839898
$ws.Range("C3").Value = 31.908378433147845
$ws.Range("D3").Value = 42.837996626649741
$ws.Range("E3").Value = 54.722123549108382

# Update the selection to match the newly edited range
$ws.Range("B1:E3").Select()
